# issue #5: stock data output to json file
#
# The stock ("股票") sheet gains a new "property_category" column (value
# "stock" for every data row), inserted between the existing "total" and
# "date" columns. Everything to the right (date, legislator_name,
# legislator_id) shifts one column over.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a blank column at H; existing H:J (date, legislator_name,
# legislator_id) shift right to I:K, inheriting their formatting.
$ws.Columns(8).Insert()

# Header
$ws.Cells.Item(1, 8).Value = "property_category"

# Data rows
$ws.Cells.Item(2, 8).Value = "stock"
$ws.Cells.Item(3, 8).Value = "stock"
